$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.683.18"
$ws.Range("E2").Value = "  +0.98%  "
# Row 3
$ws.Range("D3").Value = "1.851.18"
$ws.Range("E3").Value = "  +0.11%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.01%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "262.80"
$ws.Range("E5").Value = "  -0.52%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5382"
$ws.Range("E7").Value = "  +3.41%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3198"
$ws.Range("E8").Value = "  -2.18%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06985"
$ws.Range("E9").Value = "  +2.61%  "
# Row 10
$ws.Range("E10").Value = "  +1.55%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7754"
$ws.Range("E11").Value = "  -0.16%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07833"
$ws.Range("E12").Value = "  +0.76%  "
# Row 13
$ws.Range("D13").Value = "1.852.33"
$ws.Range("E13").Value = "  -0.29%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.62"
$ws.Range("E14").Value = "  +1.72%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.055"
$ws.Range("E15").Value = "  +0.81%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.18"
# Row 17
$ws.Range("E17").Value = "  +0.00%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008030"
$ws.Range("E18").Value = "  +0.72%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.04%  "
# Row 20
$ws.Range("D20").Value = "26.700.54"
$ws.Range("E20").Value = "  +0.97%  "
# Row 21
$ws.Range("D21").Value = "2.082.09"
$ws.Range("E21").Value = "  +0.03%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.659"
$ws.Range("E22").Value = "  +0.40%  "
# Row 23
$ws.Range("E23").Value = "  +0.96%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.421"
$ws.Range("E24").Value = "  -1.36%  "
# Row 25
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.225"
$ws.Range("E25").Value = "  +0.38%  "
# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.76"
$ws.Range("E26").Value = "  -1.49%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.701"
$ws.Range("E27").Value = "  +2.62%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.16"
$ws.Range("E28").Value = "  +0.97%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.78"
$ws.Range("E29").Value = "  -0.13%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.338"
$ws.Range("E30").Value = "  +3.57%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08769"
$ws.Range("E31").Value = "  +0.25%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.126"
$ws.Range("E32").Value = "  -0.22%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04882"
$ws.Range("E33").Value = "  +0.92%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7393"
$ws.Range("E34").Value = "  +2.92%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.144"
$ws.Range("E35").Value = "  +0.68%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.902"
$ws.Range("E36").Value = "  +1.75%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.113"
$ws.Range("E37").Value = "  +0.45%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.365"
$ws.Range("E38").Value = "  +7.82%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01753"
$ws.Range("E39").Value = "  -1.37%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4853"
$ws.Range("E40").Value = "  -0.18%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9087"
$ws.Range("E41").Value = "  -1.42%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.66"
$ws.Range("E42").Value = "  -1.18%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.924"
$ws.Range("E43").Value = "  -2.36%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.03%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.749"
$ws.Range("E45").Value = "  +0.67%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4200"
$ws.Range("E46").Value = "  +0.77%  "
# Row 47
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1255"
$ws.Range("E47").Value = "  +1.38%  "
# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.094"
$ws.Range("E48").Value = "  -0.29%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.25"
$ws.Range("E49").Value = "  +1.13%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05840"
$ws.Range("E50").Value = "  -1.59%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9008"
$ws.Range("E51").Value = "  +1.17%  "
